# Add test cases for "Become Partner" page to the KeywordDictionary sheet:
#  - elementShouldNotBePresent (Any field type)
#  - isButtonEnabled / isButtonDisabled (Button field type)
# Also adds a couple of helper methods/rows as described in the commit
# message ("Add Methods to check if button is enabled or disabled").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: elementShouldNotBePresent ------------------------------------
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Any"
$ws.Range("C19").Value = "elementShouldNotBePresent"
$ws.Range("D19").Value = "Accepts no parameters and verifies element is not available in DOM. Returns true if element is not available in DOM"

# --- Rows 20 & 21: isButtonEnabled / isButtonDisabled -----------------------
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("B20").Value = "Button"
$ws.Range("C20").Value = "isButtonEnabled"
$ws.Range("C21").Value = "isButtonDisabled"
$ws.Range("B21").Value = "Button"
$ws.Range("D20").Value = "Accepts one parameter @locator. Checks if the Button is enabled. Returns true if the button is enabled"
$ws.Range("D21").Value = "Accepts one parameter @locator. Checks if the Button is disabled Returns true if the button is disabled"

# Match formatting used by the rest of the table: column C holds the
# function name as plain Text, columns B & D wrap their text, and every
# new data row uses the same row height as its neighbours.
$ws.Range("C19").NumberFormat = "@"
$ws.Range("D19:D21").WrapText = $true

$ws.Rows.Item(19).RowHeight = 29
$ws.Rows.Item(20).RowHeight = 29
$ws.Rows.Item(21).RowHeight = 29

# Keep the visible selection in sync with where the new data ends, as
# happens naturally after typing the last new cell in Excel.
$ws.Range("D21").Select() | Out-Null
